$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# URL: hl7.fr/fhir/fr/medication -> hl7.fr/ig/fhir/medication
$meta.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-treatment-intent"

# Date updated
$meta.Range("B8").Value = "2026-01-15T08:54:26+00:00"

# Jurisdiction value was blank, now "FRANCE"
$meta.Range("B11").Value = "FRANCE"

# Description text: fix FrInpatientMedicationRequest -> FRInpatientMedicationRequest
$meta.Range("B12").Value = "Le jeu de valeurs à utiliser pour coder l'élément *treatmentIntent* de la ressource *FRInpatientMedicationRequest*."

$wb.Save()
